$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.489.31'
$ws.Range('E2').Value = '  -0.97%  '

$ws.Range('D3').Value = '1.849.75'
$ws.Range('E3').Value = '  -0.52%  '

$ws.Range('D4').Value = "'0.9985"
$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').Value = "'241.90"
$ws.Range('E5').Value = '  -1.19%  '

$ws.Range('D6').Value = "'0.6282"
$ws.Range('E6').Value = '  -2.32%  '

$ws.Range('D7').Value = "'0.9993"
$ws.Range('E7').Value = '  -0.09%  '

$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D8').Value = "'0.07539"
$ws.Range('E8').Value = '  -0.05%  '

$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').Value = "'0.2979"
$ws.Range('E9').Value = '  +0.16%  '

$ws.Range('B10').Value = 'Solana'
$ws.Range('C10').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').Value = "'24.34"
$ws.Range('E10').Value = '  -1.05%  '

$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = "'0.07696"
$ws.Range('E11').Value = '  +0.15%  '

$ws.Range('D12').Value = '1.922.85'
$ws.Range('E12').Value = '  +3.20%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'5.007"
$ws.Range('E13').Value = '  -0.83%  '

$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = "'0.6864"
$ws.Range('E14').Value = '  -0.83%  '

$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').Value = "'83.72"
$ws.Range('E15').Value = '  -0.37%  '

$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = "'0.000009812"
$ws.Range('E16').Value = '  -1.03%  '

$ws.Range('D17').Value = '2.231.06'
$ws.Range('E17').Value = '  +5.56%  '

$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').Value = "'6.233"
$ws.Range('E18').Value = '  +1.55%  '

$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '29.599.07'
$ws.Range('E19').Value = '  -0.71%  '

$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = "'234.17"
$ws.Range('E20').Value = '  -1.17%  '

$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = "'12.50"
$ws.Range('E21').Value = '  -1.29%  '

$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = "'1.001"
$ws.Range('E22').Value = '  +0.08%  '

$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = "'7.635"
$ws.Range('E23').Value = '  +1.25%  '

$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').Value = "'0.9990"
$ws.Range('E24').Value = '  -0.16%  '

$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = "'155.04"
$ws.Range('E25').Value = '  -2.41%  '

$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').Value = "'0.1393"
$ws.Range('E26').Value = '  -2.05%  '

$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = "'8.436"
$ws.Range('E27').Value = '  -1.44%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = "'17.72"
$ws.Range('E28').Value = '  -1.21%  '

$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').Value = "'1.478"
$ws.Range('E29').Value = '  -1.38%  '

$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').Value = "'0.05846"
$ws.Range('E30').Value = '  -6.72%  '

$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').Value = "'1.263"
$ws.Range('E31').Value = '  -2.22%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = "'4.099"
$ws.Range('E32').Value = '  -1.31%  '

$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = "'4.043"
$ws.Range('E33').Value = '  -1.41%  '

$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = "'1.902"
$ws.Range('E34').Value = '  -0.15%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = "'1.171"
$ws.Range('E35').Value = '  -0.20%  '

$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = "'0.7201"
$ws.Range('E36').Value = '  -1.46%  '

$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').Value = "'2.586"
$ws.Range('E37').Value = '  -0.95%  '

$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '1.240.35'
$ws.Range('E38').Value = '  +1.77%  '

$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = "'2.793"
$ws.Range('E39').Value = '  -0.94%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = "'0.01790"
$ws.Range('E40').Value = '  +0.09%  '

$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = "'0.9082"
$ws.Range('E41').Value = '  -1.32%  '

$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = "'6.149"
$ws.Range('E42').Value = '  -2.77%  '

$ws.Range('B43').Value = 'RocketPoolETH'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D43').Value = '2.121.50'
$ws.Range('E43').Value = '  +4.80%  '

$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = "'0.9992"
$ws.Range('E44').Value = '  -0.14%  '

$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = "'102.09"
$ws.Range('E45').Value = '  +0.00%  '

$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = "'67.21"
$ws.Range('E46').Value = '  -0.14%  '

$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').Value = "'7.339"
$ws.Range('E47').Value = '  +9.03%  '

$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = "'1.730"
$ws.Range('E48').Value = '  +3.07%  '

$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').Value = "'0.4037"
$ws.Range('E49').Value = '  -0.65%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'9.142"
$ws.Range('E50').Value = '  -0.87%  '

$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = "'0.00000000117"
$ws.Range('E51').Value = '  -1.40%  '
